$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(42606.571250000001, 16, 75, 23, 25, 75, 5609, 2134, 427, 73, 23, 2, 6, "Noun")

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(3, $col).Value = $values[$i]
}
